$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-44 and add new rows 45-61 with refreshed Entsoe consumption data
$data = @(
    @(2, 5945, 45981.95833333334),
    @(3, 5897, 45981.96875),
    @(4, 5855, 45981.97916666666),
    @(5, 5806, 45981.98958333334),
    @(6, 5758, 45982),
    @(7, 5705, 45982.01041666666),
    @(8, 5655, 45982.02083333334),
    @(9, 5652, 45982.03125),
    @(10, 5625, 45982.04166666666),
    @(11, 5586, 45982.05208333334),
    @(12, 5576, 45982.0625),
    @(13, 5578, 45982.07291666666),
    @(14, 5552, 45982.08333333334),
    @(15, 5546, 45982.10416666666),
    @(16, 5577, 45982.11458333334),
    @(17, 5583, 45982.125),
    @(18, 5624, 45982.13541666666),
    @(19, 5644, 45982.14583333334),
    @(20, 5687, 45982.15625),
    @(21, 5715, 45982.16666666666),
    @(22, 5775, 45982.17708333334),
    @(23, 5860, 45982.1875),
    @(24, 5961, 45982.19791666666),
    @(25, 6188, 45982.20833333334),
    @(26, 6326, 45982.21875),
    @(27, 6458, 45982.22916666666),
    @(28, 6656, 45982.23958333334),
    @(29, 6906, 45982.25),
    @(30, 7043, 45982.26041666666),
    @(31, 7201, 45982.27083333334),
    @(32, 7228, 45982.28125),
    @(33, 7346, 45982.29166666666),
    @(34, 7393, 45982.30208333334),
    @(35, 7380, 45982.3125),
    @(36, 7376, 45982.32291666666),
    @(37, 7339, 45982.33333333334),
    @(38, 7327, 45982.34375),
    @(39, 7276, 45982.35416666666),
    @(40, 7293, 45982.36458333334),
    @(41, 7202, 45982.375),
    @(42, 7148, 45982.38541666666),
    @(43, 7048, 45982.39583333334),
    @(44, 7027, 45982.40625),
    @(45, 6902, 45982.41666666666),
    @(46, 6894, 45982.42708333334),
    @(47, 6858, 45982.4375),
    @(48, 6760, 45982.44791666666),
    @(49, 6756, 45982.45833333334),
    @(50, 6718, 45982.46875),
    @(51, 6805, 45982.47916666666),
    @(52, 6869, 45982.48958333334),
    @(53, 6875, 45982.5),
    @(54, 6932, 45982.51041666666),
    @(55, 6909, 45982.52083333334),
    @(56, 6899, 45982.53125),
    @(57, 6955, 45982.54166666666),
    @(58, 6965, 45982.55208333334),
    @(59, 6978, 45982.5625),
    @(60, 6964, 45982.57291666666),
    @(61, 6949, 45982.58333333334)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

